$p = $ppt.ActivePresentation

# --- 1. Remove trailing slides (presentation slides 16 through 33) ---
for ($i = $p.Slides.Count; $i -ge 16; $i--) {
    $p.Slides.Item($i).Delete()
}

# --- 2. Slide 14: "Jump to document" edits ---
$s14 = $p.Slides.Item(14)
$shape14 = $s14.Shapes.Item(2)
$tf14 = $shape14.TextFrame
$tr14 = $tf14.TextRange

$tr14.Paragraphs(1,1).Text = "Jump to document.`r"
$tr14.InsertAfter("`r")
$tr14.InsertAfter("Start with table of model runs.`r")

# --- 3. Slide 15: add "Fix stage-2 selectivity..." bullet ---
$s15 = $p.Slides.Item(15)
$shape15 = $s15.Shapes.Item(2)
$tf15 = $shape15.TextFrame
$tr15 = $tf15.TextRange

$tr15.Paragraphs(2,1).InsertAfter("Fix stage-2 selectivity at 1.0 (rather than stage-1) if estimated stage-2 selectivity > 1.0`r")
